$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column N formatting to column O for rows 3-10
$ws.Range("N3:N10").Copy()
$ws.Range("O3:O10").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Set values for new column O
$ws.Range("O4").Value = 2021
$ws.Range("O6").Value = 1860
$ws.Range("O7").Value = 1
$ws.Range("O8").Value = 510
$ws.Range("O9").Value = 178
$ws.Range("O10").Value = 821

# Update selection to match target
$ws.Range("P9").Select()
